$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so values like "1.000" or "1.001" are not
# reinterpreted as numbers; ClearFormats afterwards restores the original (default)
# cell style so no spurious style/number-format metadata is introduced.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.237.59"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.856.97"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "0.7062"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").Value = "237.83"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.07997"
$ws.Range("E8").Value = "  +2.66%  "
$ws.Range("D9").Value = "0.3018"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").Value = "23.47"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").Value = "0.08171"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "1.870.98"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "5.186"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "0.7043"
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").Value = "89.60"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "29.244.75"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "0.000007927"
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").Value = "5.792"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "13.21"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").Value = "238.08"
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("D21").Value = "0.9982"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "2.101.90"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "7.472"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("D25").Value = "162.92"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").Value = "8.885"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").Value = "0.1428"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").Value = "18.07"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "1.918"
$ws.Range("E29").Value = "  -2.92%  "
$ws.Range("D30").Value = "1.428"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("D31").Value = "1.477"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "4.361"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").Value = "4.016"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "0.05187"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Value = "1.157"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("D36").Value = "0.7126"
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  -2.28%  "
$ws.Range("D38").Value = "2.653"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "0.01850"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").Value = "2.722"
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("D41").Value = "0.9408"
$ws.Range("E41").Value = "  +2.96%  "
$ws.Range("D42").Value = "1.134.52"
$ws.Range("E42").Value = "  +4.32%  "
$ws.Range("D43").Value = "5.941"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").Value = "0.4253"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "70.29"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "102.90"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").Value = "0.5306"
$ws.Range("E48").Value = "  -4.25%  "
$ws.Range("D49").Value = "1.760"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").Value = "2.013.19"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").Value = "9.162"
$ws.Range("E51").Value = "  +0.09%  "

$dRange.ClearFormats()
